$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.847.88'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.45%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.634.57'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.23%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.78%  '

$ws.Range("E7").Value = '  +0.38%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  -1.72%  '

$ws.Range("E10").Value = '  +1.40%  '

$ws.Range("E11").Value = '  +0.32%  '

$ws.Range("E12").Value = '  +0.95%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.67'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.109.90'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.09%  '

$ws.Range("E15").Value = '  +0.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.750.09'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.50%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.634.31'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.75%  '

$ws.Range("E18").Value = '  +0.91%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.52'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.92%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '343.74'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.62%  '

$ws.Range("E22").Value = '  +0.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.24%  '

$ws.Range("E24").Value = '  +8.42%  '

$ws.Range("E25").Value = '  +5.11%  '

$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.81%  '

$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.26'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '581.42'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.83%  '

$ws.Range("E29").Value = '  +5.15%  '

$ws.Range("E30").Value = '  +1.00%  '

$ws.Range("E31").Value = '  -0.21%  '

$ws.Range("E32").Value = '  +0.27%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.73'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.34%  '

$ws.Range("E34").Value = '  +2.93%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.44'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.68%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.403'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.04%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.80'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.00%  '

$ws.Range("E39").Value = '  +3.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '153.41'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.27%  '

$ws.Range("E41").Value = '  +8.89%  '

$ws.Range("E42").Value = '  +0.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '162.27'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.29%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '24.15'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.85%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.93'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.17%  '

$ws.Range("E46").Value = '  -0.18%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.633'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.21%  '

$ws.Range("E48").Value = '  -1.21%  '

$ws.Range("E49").Value = '  -0.58%  '

$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0237'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.69%  '

$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.790'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.49%  '
